$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells retain their text formatting (avoid Excel auto-converting
# numeric-looking strings like "1.00" or "62.692.71" into numbers).
$textCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "B26", "C26", "D26", "E26", "B27", "C27", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "B33", "C33", "D33", "E33", "B34", "C34", "D34", "E34", "B35", "C35", "D35", "E35", "D36", "E36", "D37", "E37", "B38", "C38", "D38", "E38", "B39", "C39", "D39", "E39", "D40", "E40", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "B46", "C46", "D46", "E46", "B47", "C47", "D47", "E47", "D48", "E48", "D50", "E50", "D51", "E51"
)
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '62.692.71'
$ws.Range('E2').Value = '  +1.42%  '
$ws.Range('D3').Value = '2.922.96'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = '570.49'
$ws.Range('E5').Value = '  -2.77%  '
$ws.Range('D6').Value = '145.57'
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '2.921.27'
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').Value = '0.503'
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').Value = '6.96'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').Value = '0.151'
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').Value = '0.435'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').Value = '0.0000241'
$ws.Range('E13').Value = '  +1.43%  '
$ws.Range('D14').Value = '32.90'
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').Value = '3.390.95'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '62.410.50'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').Value = '6.63'
$ws.Range('E18').Value = '  +0.60%  '
$ws.Range('D19').Value = '2.893.74'
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('D20').Value = '433.34'
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('D21').Value = '13.14'
$ws.Range('E21').Value = '  -2.12%  '
$ws.Range('D22').Value = '0.664'
$ws.Range('E22').Value = '  +0.64%  '
$ws.Range('D23').Value = '6.90'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').Value = '78.94'
$ws.Range('E24').Value = '  -2.39%  '
$ws.Range('D25').Value = '12.08'
$ws.Range('E25').Value = '  +1.31%  '
$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D26').Value = '10.05'
$ws.Range('E26').Value = '  -1.66%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').Value = '2.03'
$ws.Range('E28').Value = '  -2.00%  '
$ws.Range('D29').Value = '0.0000113'
$ws.Range('E29').Value = '  +5.42%  '
$ws.Range('D30').Value = '7.02'
$ws.Range('E30').Value = '  -2.68%  '
$ws.Range('D31').Value = '2.50'
$ws.Range('E31').Value = '  -2.27%  '
$ws.Range('D32').Value = '2.03'
$ws.Range('E32').Value = '  -3.51%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.107'
$ws.Range('E33').Value = '  -1.70%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '25.90'
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').Value = '0.996'
$ws.Range('E35').Value = '  -0.46%  '
$ws.Range('D36').Value = '0.954'
$ws.Range('E36').Value = '  -2.27%  '
$ws.Range('D37').Value = '5.43'
$ws.Range('E37').Value = '  -1.40%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').Value = '2.98'
$ws.Range('E38').Value = '  -2.83%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '49.05'
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('D40').Value = '1.91'
$ws.Range('E40').Value = '  -4.66%  '
$ws.Range('E41').Value = '  -0.88%  '
$ws.Range('D42').Value = '41.26'
$ws.Range('E42').Value = '  +6.16%  '
$ws.Range('D43').Value = '8.12'
$ws.Range('E43').Value = '  -2.77%  '
$ws.Range('D44').Value = '0.269'
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('D45').Value = '2.708.85'
$ws.Range('E45').Value = '  +0.41%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0339'
$ws.Range('E46').Value = '  +0.87%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '133.23'
$ws.Range('E47').Value = '  -1.18%  '
$ws.Range('D48').Value = '351.26'
$ws.Range('E48').Value = '  +2.13%  '
$ws.Range('D50').Value = '0.000220'
$ws.Range('E50').Value = '  +14.96%  '
$ws.Range('D51').Value = '0.103'
$ws.Range('E51').Value = '  -0.72%  '
